$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.020.86"
$ws.Range("E2").Value = "  +3.76%  "

$ws.Range("D3").Value = "1.692.53"
$ws.Range("E3").Value = "  +3.58%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.534"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.68%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.42"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.50%  "

$ws.Range("E9").Value = "  +2.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0641"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.17%  "

$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("D12").Value = "1.935.62"
$ws.Range("E12").Value = "  +3.59%  "

$ws.Range("D13").Value = "1.697.45"
$ws.Range("E13").Value = "  +3.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.608"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.33%  "

$ws.Range("D17").Value = "31.039.11"
$ws.Range("E17").Value = "  +3.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.91%  "

$ws.Range("D20").Value = "0.0₃0721"
$ws.Range("E20").Value = "  +2.49%  "

$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("E22").Value = "  +3.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.36%  "

$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.44%  "

$ws.Range("E26").Value = "  +2.76%  "

$ws.Range("E27").Value = "  +2.44%  "

$ws.Range("E28").Value = "  +1.61%  "

$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0502"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.42%  "

$ws.Range("E32").Value = "  +3.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.04%  "

$ws.Range("D34").Value = "1.515.98"
$ws.Range("E34").Value = "  +6.59%  "

$ws.Range("E35").Value = "  +2.91%  "

$ws.Range("E36").Value = "  +0.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "83.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.58%  "

$ws.Range("E38").Value = "  +10.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0180"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.41%  "

$ws.Range("E40").Value = "  -4.07%  "

$ws.Range("E41").Value = "  +0.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.847"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.72%  "

$ws.Range("E44").Value = "  +0.69%  "

$ws.Range("E45").Value = "  +2.88%  "

$ws.Range("E46").Value = "  +0.15%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.42%  "

$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("E48").Value = "  +6.14%  "

$ws.Range("D49").Value = "1.825.13"
$ws.Range("E49").Value = "  +2.74%  "

$ws.Range("E50").Value = "  +9.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "93.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.95%  "
